$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New record goes right after the last existing data row (row 73).
$row = 74

$ws.Cells.Item($row, 1).Value = "Kindergarden"
$ws.Cells.Item($row, 2).Value = "Kindergarden Bussum Kamerlingh Onnesweg"
$ws.Cells.Item($row, 3).Value = "KDV"
$ws.Cells.Item($row, 4).Value = "'2024-09-03"
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
